$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

function Set-TextCell($addr, $value) {
    # Numeric-looking text must be forced to stay text (leading apostrophe,
    # matching how Excel preserves user-typed text that looks like a number).
    $ws.Range($addr).Value = "'" + $value
}

# Row 2
Set-Cell "D2" "66.416.74"
Set-Cell "E2" "  -0.84%  "

# Row 3
Set-Cell "D3" "3.464.29"
Set-Cell "E3" "  -1.74%  "

# Row 4
Set-Cell "E4" "  +0.01%  "

# Row 5
Set-TextCell "D5" "586.30"
Set-Cell "E5" "  -0.07%  "

# Row 6
Set-TextCell "D6" "176.63"
Set-Cell "E6" "  -0.88%  "

# Row 7
Set-TextCell "D7" "0.625"
Set-Cell "E7" "  +3.45%  "

# Row 8
Set-Cell "E8" "  +0.02%  "

# Row 9
Set-Cell "D9" "3.460.82"
Set-Cell "E9" "  -1.82%  "

# Row 10
Set-TextCell "D10" "0.133"
Set-Cell "E10" "  -1.56%  "

# Row 11
Set-TextCell "D11" "6.97"
Set-Cell "E11" "  +0.51%  "

# Row 12
Set-Cell "E12" "  -2.01%  "

# Row 13
Set-Cell "D13" "4.068.69"
Set-Cell "E13" "  -1.75%  "

# Row 14
Set-Cell "B14" "TRON"
Set-Cell "C14" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D14" "0.134"
Set-Cell "E14" "  +1.24%  "

# Row 15
Set-Cell "B15" "Avalanche"
Set-Cell "C15" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D15" "30.35"
Set-Cell "E15" "  -1.20%  "

# Row 16
Set-Cell "D16" "66.312.15"
Set-Cell "E16" "  -0.95%  "

# Row 17
Set-Cell "E17" "  -1.03%  "

# Row 18
Set-Cell "D18" "3.471.43"
Set-Cell "E18" "  -1.55%  "

# Row 19
Set-Cell "E19" "  -1.81%  "

# Row 20
Set-TextCell "D20" "13.82"
Set-Cell "E20" "  -1.96%  "

# Row 21
Set-TextCell "D21" "373.13"
Set-Cell "E21" "  -2.64%  "

# Row 22
Set-Cell "E22" "  -2.49%  "

# Row 23
Set-TextCell "D23" "73.27"
Set-Cell "E23" "  +1.03%  "

# Row 24
Set-TextCell "D24" "0.999"
Set-Cell "E24" "  -0.07%  "

# Row 25
Set-TextCell "D25" "0.0000126"
Set-Cell "E25" "  +4.07%  "

# Row 26
Set-TextCell "D26" "0.536"
Set-Cell "E26" "  -3.06%  "

# Row 27
Set-Cell "D27" "3.618.23"
Set-Cell "E27" "  -1.36%  "

# Row 28
Set-TextCell "D28" "9.94"
Set-Cell "E28" "  +0.32%  "

# Row 29
Set-Cell "E29" "  +2.86%  "

# Row 30
Set-TextCell "D30" "1.00"
Set-Cell "E30" "  +0.00%  "

# Row 31
Set-Cell "E31" "  -0.51%  "

# Row 32
Set-Cell "E32" "  -1.09%  "

# Row 33
Set-TextCell "D33" "23.73"
Set-Cell "E33" "  -4.19%  "

# Row 34
Set-TextCell "D34" "0.999"
Set-Cell "E34" "  -0.03%  "

# Row 35
Set-Cell "E35" "  -3.40%  "

# Row 36
Set-Cell "E36" "  -5.81%  "

# Row 37
Set-TextCell "D37" "1.55"
Set-Cell "E37" "  -2.44%  "

# Row 38
Set-TextCell "D38" "161.18"

# Row 39
Set-TextCell "D39" "0.886"
Set-Cell "E39" "  -1.47%  "

# Row 40
Set-TextCell "D40" "28.30"
Set-Cell "E40" "  -6.79%  "

# Row 41
Set-TextCell "D41" "1.81"
Set-Cell "E41" "  +0.02%  "

# Row 42
Set-Cell "B42" "Maker"
Set-Cell "C42" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-Cell "D42" "2.784.72"
Set-Cell "E42" "  +1.68%  "

# Row 43
Set-Cell "B43" "Filecoin"
Set-Cell "C43" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D43" "4.51"
Set-Cell "E43" "  -0.85%  "

# Row 44
Set-Cell "B44" "dogwifhat"
Set-Cell "C44" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D44" "2.58"
Set-Cell "E44" "  +0.28%  "

# Row 45
Set-TextCell "D45" "6.46"
Set-Cell "E45" "  -2.99%  "

# Row 46
Set-Cell "E46" "  -1.93%  "

# Row 47
Set-TextCell "D47" "25.33"
Set-Cell "E47" "  -0.17%  "

# Row 48
Set-TextCell "D48" "340.32"
Set-Cell "E48" "  +5.00%  "

# Row 49
Set-TextCell "D49" "40.03"
Set-Cell "E49" "  -1.97%  "

# Row 50
Set-TextCell "D50" "0.0293"
Set-Cell "E50" "  -1.97%  "

# Row 51
Set-Cell "E51" "  +0.58%  "

